# Oral defence 07-01-2021.pptx
# Slide 13 ("Same approximation of p, but decision maker is more uncertain
# than the observer") textbox "Tekstfelt 13": split the leading run
# "Same approximation of p, but " into three runs so that "of p," becomes
# "of x," while "Same approximation " and "but " keep their own runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Tekstfelt 13") {
        $target = $sh
        break
    }
}
if ($target -eq $null) {
    $target = $s.Shapes.Item(6)
}

$tr = $target.TextFrame.TextRange

# "Same approximation of p, but decision maker is more uncertain than the observer"
#  1-19 = "Same approximation "
# 20-25 = "of p, "
# 26-29 = "but "
$tr.Characters(20, 6).Text = "of x, "
$tr.Characters(26, 4).Text = "but "
